# Add a "Valor Total" summary row below the existing data table,
# formatted the same way as the header row (A1:D1), and merge it
# across columns A:D just like the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row goes right after the last existing data row (row 10).
$targetRow = 11

# Merge the destination range first...
$ws.Range("A$targetRow`:D$targetRow").Merge()

# ...then copy the formatting (styles/borders/font/alignment) from the
# title row (A1:D1) onto the new row, reusing the existing cell styles
# instead of creating new ones.
$ws.Range("A1:D1").Copy()
$ws.Range("A$targetRow").PasteSpecial(-4122)

# Finally set the text for the summary cell.
$ws.Range("A$targetRow").Value = "Valor Total: 299"

Write-Host "Added summary row $targetRow with total value"
